# Update for last Shaconemo revision - same variables as previously this time
#
# The shared-string note "Identified in one of the shaconemo (238) ping
# files. " is used by every row that was matched against that ping-file
# revision. Bump the revision number from 238 to 239 everywhere it is used
# (this is exactly what Excel's Find&Replace "Replace All" does - it updates
# the single shared string instance, which is reflected in every cell that
# references it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "Identified in one of the shaconemo (238) ping files. "
$newText = "Identified in one of the shaconemo (239) ping files. "

$ws.Cells.Replace($oldText, $newText)

# Mirror the post-"Replace All" UI state: Excel leaves the selection sitting
# on the matched cells, with the active cell on the first match (G13, the
# "Sea Surface Total Chlorophyll Mass Concentration" row).
$ws.Range("G13").Select()
